$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Range("A1").End(4).Row
if ($lastRow -lt 2) { $lastRow = 51 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    if ($cell.Text -eq "2025-03-02 00:17:04") {
        $cell.Value = "2025-03-02 00:30:04"
    }
}
